$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip leading tab character from text in B1:C5
for ($r = 1; $r -le 5; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    $bCell.Value = $bVal.TrimStart("`t")

    $cCell = $ws.Cells.Item($r, 3)
    $cVal = $cCell.Value()
    $cCell.Value = $cVal.TrimStart("`t")
}

# Delete row 6 entirely
$ws.Rows.Item(6).Delete()
